# Updated vm_pu results for Case_3_241 (380 kV case) -- res_bus/vm_pu sheet
# Row key = spreadsheet row number; inner hashtable maps column letter -> new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ B=1.02; C=1.031113416865132; D=1.032837360698306; E=1.040277374940922; F=1.050289666076885; I=1.031449058011148; J=1.03625093380873; K=1.035641440296128; L=1.043060179436843; M=1.053044355885616; N=1.016022855641347 }
    3 = @{ B=1.02; C=1.032332122473596; D=1.033714814677224; E=1.041432717982604; F=1.051700288617853; I=1.031686677949418; J=1.037110006337224; K=1.036327566236724; L=1.044025005295094; M=1.054265853680189; N=1.016311031971012 }
    4 = @{ B=1.02; C=1.033119910418963; D=1.034281474536655; E=1.042179999297862; F=1.052613168698612; I=1.031838173716436; J=1.037664607699316; K=1.036769792992018; L=1.044648435221412; M=1.055055833857171; N=1.016496967585057 }
    5 = @{ B=1.02; C=1.033450908599552; D=1.034519433200118; E=1.04249408595799; F=1.052996973626629; I=1.031901321724081; J=1.037897458234611; K=1.036955288944549; L=1.044910317129041; M=1.05538784586918; N=1.016575007611646 }
    6 = @{ B=1.02; C=1.033506473712156; D=1.034559371986415; E=1.04254681842271; F=1.053061418020188; I=1.031911892847079; J=1.037936537090846; K=1.036986410155474; L=1.044954276109564; M=1.055443586614147; N=1.016588103423907 }
    7 = @{ B=1.02; C=1.033124333967708; D=1.034284655192632; E=1.042184196413229; F=1.052618296996567; I=1.031839019627901; J=1.037667720251635; K=1.036772273232087; L=1.044651935315366; M=1.055060270590519; N=1.016498010859356 }
    8 = @{ B=1.02; C=1.031525450600251; D=1.03313413111552; E=1.040667892518585; F=1.050766371045253; I=1.031529830953437; J=1.036541526566387; K=1.035873681200929; L=1.043386429574378; M=1.053457253388382; N=1.016120356982083 }
    9 = @{ B=1.02; C=1.028701794608885; D=1.031098203195004; E=1.037993556363599; F=1.047503766818822; I=1.030967676582242; J=1.03454718403276; K=1.034276854618077; L=1.04114964458908; M=1.050629268905019; N=1.015450773251903 }
    10 = @{ B=1.02; C=1.026815001981489; D=1.029735098402226; E=1.036208905852318; F=1.045328987945316; I=1.030581237374241; J=1.033210898311329; K=1.033203231086761; L=1.039653758196734; M=1.048741582956589; N=1.01500159001646 }
    11 = @{ B=1.02; C=1.025996927335122; D=1.029143462835032; E=1.035435683912792; F=1.044387308862018; I=1.030411131501092; J=1.032630652816476; K=1.032736172226834; L=1.039004880743268; M=1.04792359524019; N=1.014806418720553 }
    12 = @{ B=1.02; C=1.025692892058116; D=1.028923490883996; E=1.035148403483291; F=1.04403752555173; I=1.030347529097152; J=1.032414877478747; K=1.032562357820991; L=1.038763683642784; M=1.047619663541212; N=1.014733821783496 }
    13 = @{ B=1.02; C=1.025758116222522; D=1.028970685239066; E=1.035210029385436; F=1.044112555402961; I=1.030361190931203; J=1.032461173203332; K=1.032599656465975; L=1.038815429205145; M=1.04768486225998; N=1.014749398682605 }
    14 = @{ B=1.02; C=1.025971799089628; D=1.029125284221638; E=1.035411938695081; F=1.044358395723528; I=1.030405882626961; J=1.032612821797472; K=1.032721811372717; L=1.038984946902871; M=1.04789847410156; N=1.014800419912155 }
    15 = @{ B=1.02; C=1.026103434063577; D=1.029220509599027; E=1.035536332143108; F=1.044509865705588; I=1.030433363314377; J=1.032706224794886; K=1.032797031540299; L=1.039089369069661; M=1.048030074764182; N=1.014831842288818 }
    16 = @{ B=1.02; C=1.026869271826195; D=1.029774333636567; E=1.036260212190391; F=1.045391484053825; I=1.030592468222937; J=1.033249372870992; K=1.033234182338281; L=1.039696797637182; M=1.04879585702293; N=1.015014528677577 }
    17 = @{ B=1.02; C=1.027349369730185; D=1.030121356366784; E=1.036714158839086; F=1.045944501055267; I=1.030691527141455; J=1.033589638718404; K=1.033507812705638; L=1.040077512187288; M=1.049276047025031; N=1.015128942641409 }
    18 = @{ B=1.02; C=1.027629298655574; D=1.030323633665831; E=1.036978894384994; F=1.046267068186128; I=1.030749038811454; J=1.033787953395551; K=1.033667207131372; L=1.040299465876172; M=1.049556075371252; N=1.015195613579289 }
    19 = @{ B=1.02; C=1.027724729640601; D=1.030392582111627; E=1.037069155010031; F=1.046377055643927; I=1.030768603411467; J=1.033855547073244; K=1.033721521015211; L=1.040375127653179; M=1.049651548004602; N=1.015218335676403 }
    20 = @{ B=1.02; C=1.02729787055883; D=1.030084138106909; E=1.036665459217177; F=1.045885167421964; I=1.030680926753227; J=1.033553147631532; K=1.033478476442102; L=1.040036676585284; M=1.049224533244439; N=1.015116673812747 }
    21 = @{ B=1.02; C=1.025908879420014; D=1.029079764543283; E=1.035352483428616; F=1.044286001957252; I=1.030392733566431; J=1.032568171883743; K=1.032685848862572; L=1.038935033020692; M=1.04783557337836; N=1.014785398242549 }
    22 = @{ B=1.019999999999999; C=1.025034603790201; D=1.028447046779172; E=1.034526549097657; F=1.043280530942875; I=1.03020911920151; J=1.031947452274098; K=1.03218559426623; L=1.038241371512682; M=1.046961729921465; N=1.014576523820779 }
    23 = @{ B=1.02; C=1.025498166190807; D=1.028782579365439; E=1.03496443272582; F=1.043813552635479; I=1.030306685842244; J=1.032276643414703; K=1.032450969128791; L=1.038609191613887; M=1.047425024079205; N=1.014687308095011 }
    24 = @{ B=1.02; C=1.027321141139022; D=1.030100955854263; E=1.03668746461552; F=1.045911977731572; I=1.030685717439619; J=1.033569636866888; K=1.033491732884276; L=1.040055128777967; M=1.049247810284712; N=1.015122217768508 }
    25 = @{ B=1.02; C=1.029432532199892; D=1.031625560180302; E=1.038685236491995; F=1.048347162453626; I=1.031115062177127; J=1.035063946808031; K=1.034691266592681; L=1.041728725819444; M=1.051360774726659; N=1.015624366864009 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
